$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.239.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.233.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.57%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.92%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.98%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'74.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.47%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -3.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'42.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.71%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0961"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'6.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.56%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.29%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.569.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.28%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'14.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.24%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -2.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.275.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.099.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.52%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +3.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.43%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.22%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'11.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.58%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'230.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.42%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -6.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.02%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -2.90%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.16%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.72%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'167.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.19%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'20.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -7.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'30.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.58%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.94%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'13.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.83%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.09%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'65.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.36%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.51%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'104.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.17%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.69%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -2.73%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.76%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.58%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.441.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.49%  "
$ws.Range("E51").Style = "Normal"
